# This script updates the "UK" worksheet (second sheet) of the workbook to
# reflect a refreshed / re-run set of regression coefficients (RMSE values),
# trims the sheet from 27 data rows down to 18, widens column A, and
# clears the previous selection/active cell.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("UK")
$ws.Activate()

# New data: (row, label, coefficient)
$rows = @(
    @(2,  "Wages_FemalesNE", 0.5511963385465366),
    @(3,  "Wages_MalesNE",   0.56378226458419012),
    @(4,  "Wages_FemalesE",  0.3588496779494994),
    @(5,  "Wages_MalesE",    0.37499140584966989),
    @(6,  "I3a",             2.4158001968775116),
    @(7,  "I3b",             1.8361431734912119),
    @(8,  "I4b",             1.5969934128957666),
    @(9,  "I5b",             884.81552442493353),
    @(10, "C1b",             1.528311),
    @(11, "S1b",             1.1671),
    @(12, "S2g",             1.2093),
    @(13, "S2h",             0.98887999999999998),
    @(14, "S2i",             0.95128999999999997),
    @(15, "S2j",             0.84714999999999996),
    @(16, "S2k",             0.94330000000000003),
    @(17, "S3e",             1.2788999999999999),
    @(18, "HM1",             4.4850000000000003)
)

foreach ($r in $rows) {
    $rowIdx = $r[0]
    $label = $r[1]
    $value = $r[2]
    $ws.Cells.Item($rowIdx, 1).Value = $label
    $ws.Cells.Item($rowIdx, 2).Value = $value
}

# Rows 10-17 (C1b..S3e) got re-pasted with an explicit "No Fill" style applied.
$fillRange = $ws.Range("A10:B17")
$fillRange.Interior.Pattern = 1          # xlSolid
$fillRange.Interior.ColorIndex = -4105   # xlAutomatic (renders as "no fill")

# Remove the now-unused trailing rows (previously rows 19-27).
$ws.Rows("19:27").Delete()

# Column A was widened (22.14 chars maps to a stored width of 23).
$ws.Columns("A").ColumnWidth = 22.14

# Update the selection to match the saved view (active cell A7, no multi-select).
$ws.Range("A7").Select()

$wb.Save()
